# Apply the latest coinranking.com snapshot values to the cryptos sheet.
# (Updated cryptos list refresh - GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.263.24'
$ws.Range('E2').Value = '  +0.34%  '

$ws.Range('D3').Value = '1.907.55'
$ws.Range('E3').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.69'

$ws.Range('E6').Value = '  +0.04%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5256'
$ws.Range('E7').Value = '  +0.30%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').Value = '  +1.30%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07304'
$ws.Range('E9').Value = '  +0.75%  '

$ws.Range('E10').Value = '  +2.12%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9056'

$ws.Range('E12').Value = '  -4.06%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '96.58'
$ws.Range('E13').Value = '  +1.08%  '

$ws.Range('E14').Value = '  +1.51%  '

$ws.Range('D15').Value = '1.676.21'
$ws.Range('E15').Value = '  -12.13%  '

$ws.Range('E16').Value = '  +0.00%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008685'
$ws.Range('E17').Value = '  +0.58%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.74'
$ws.Range('E18').Value = '  +1.18%  '

$ws.Range('E19').Value = '  +0.04%  '

$ws.Range('D20').Value = '27.287.01'
$ws.Range('E20').Value = '  +0.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.124'
$ws.Range('E21').Value = '  +1.11%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.515'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.345'
$ws.Range('E24').Value = '  +2.21%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.86'
$ws.Range('E25').Value = '  +1.60%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.25'
$ws.Range('E26').Value = '  -0.05%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '116.89'
$ws.Range('E28').Value = '  +1.65%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.852'
$ws.Range('E29').Value = '  +0.65%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.872'
$ws.Range('E30').Value = '  -0.86%  '

$ws.Range('E31').Value = '  -0.27%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8244'
$ws.Range('E32').Value = '  +2.26%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05083'
$ws.Range('E33').Value = '  +0.65%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.229'
$ws.Range('E34').Value = '  -0.89%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.990'
$ws.Range('E35').Value = '  +1.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.758'
$ws.Range('E36').Value = '  +5.14%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.361'
$ws.Range('E37').Value = '  -2.54%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5766'
$ws.Range('E38').Value = '  +1.03%  '

$ws.Range('E39').Value = '  +0.34%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.084'
$ws.Range('E40').Value = '  +0.98%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.114'
$ws.Range('E41').Value = '  -0.70%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.618'
$ws.Range('E42').Value = '  -0.37%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '117.01'
$ws.Range('E43').Value = '  +0.65%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1523'
$ws.Range('E44').Value = '  +0.23%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4915'
$ws.Range('E45').Value = '  +0.96%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.16'
$ws.Range('E47').Value = '  +0.12%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.645'
$ws.Range('E48').Value = '  +1.87%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '38.77'
$ws.Range('E49').Value = '  +3.14%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06056'
$ws.Range('E50').Value = '  +1.90%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.02'
$ws.Range('E51').Value = '  -0.42%  '
